$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 2418.625
$ws.Range("I5").Value = 345.36365
$ws.Range("J5").Value = 6979.8
$ws.Range("K5").Value = 345.36365
$ws.Range("L5").Value = 6979.8
$ws.Range("M5").Value = -230.36365
$ws.Range("N5").Value = -7209.8
$ws.Range("H74").Value = 10215928
$ws.Range("I74").Value = 14294299
$ws.Range("K74").Value = 14294299
$ws.Range("M74").Value = -14293363
$ws.Range("H77").Value = 10215928
$ws.Range("I77").Value = 14294299
$ws.Range("K77").Value = 71471495
$ws.Range("M77").Value = -71466815
$ws.Range("H80").Value = 626.9375
$ws.Range("I80").Value = 805
$ws.Range("K80").Value = 2415
$ws.Range("M80").Value = -1417
$ws.Range("H83").Value = 626.9375
$ws.Range("I83").Value = 805
$ws.Range("K83").Value = 7245
$ws.Range("M83").Value = -2253
$ws.Range("H112").Value = 4656.2666
$ws.Range("J112").Value = 4526.4614
$ws.Range("L112").Value = 13579.3842
$ws.Range("N112").Value = -15795.3842
$ws.Range("H127").Value = 4877.8
$ws.Range("I127").Value = 4877.8
$ws.Range("K127").Value = 14633.4
$ws.Range("M127").Value = -9673.400000000001
$ws.Range("H132").Value = 289632.66
$ws.Range("I132").Value = 333020.06
$ws.Range("K132").Value = 999060.1799999999
$ws.Range("M132").Value = -996530.1799999999
$ws.Range("H138").Value = 4425.375
$ws.Range("I138").Value = 2166.15
$ws.Range("K138").Value = 6498.450000000001
$ws.Range("M138").Value = -1358.450000000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 935.4
$ws.Range("I4").Value = 935.4
$ws.Range("K4").Value = 935.4
$ws.Range("M4").Value = -819.4
$ws.Range("H13").Value = 21571858
$ws.Range("I13").Value = 37750000
$ws.Range("J13").Value = 1000
$ws.Range("K13").Value = 37750000
$ws.Range("L13").Value = 1000
$ws.Range("M13").Value = -37749856
$ws.Range("N13").Value = -1288
$ws.Range("H61").Value = 11956.571
$ws.Range("I61").Value = 18065.666
$ws.Range("K61").Value = 18065.666
$ws.Range("M61").Value = -17853.666
$ws.Range("H74").Value = 71433290
$ws.Range("I74").Value = 200005200
$ws.Range("K74").Value = 200005200
$ws.Range("M74").Value = -200004326
$ws.Range("H77").Value = 71433290
$ws.Range("I77").Value = 200005200
$ws.Range("K77").Value = 1000026000
$ws.Range("M77").Value = -1000021632
$ws.Range("H132").Value = 955780.7
$ws.Range("I132").Value = 1364591
$ws.Range("K132").Value = 4093773
$ws.Range("M132").Value = -4091243
$ws.Range("H136").Value = 11956.571
$ws.Range("I136").Value = 18065.666
$ws.Range("K136").Value = 54196.99800000001
$ws.Range("M136").Value = -51646.99800000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1430.3125
$ws.Range("I22").Value = 1363.7273
$ws.Range("K22").Value = 1363.7273
$ws.Range("M22").Value = -1190.7273
$ws.Range("H64").Value = 4033.6
$ws.Range("I64").Value = 914.5
$ws.Range("J64").Value = 6113
$ws.Range("K64").Value = 914.5
$ws.Range("L64").Value = 6113
$ws.Range("M64").Value = -689.5
$ws.Range("N64").Value = -6563
$ws.Range("H67").Value = 4033.6
$ws.Range("I67").Value = 914.5
$ws.Range("J67").Value = 6113
$ws.Range("K67").Value = 914.5
$ws.Range("L67").Value = 6113
$ws.Range("M67").Value = -134.5
$ws.Range("N67").Value = -7673
$ws.Range("H86").Value = 6394.5835
$ws.Range("I86").Value = 6217.3335
$ws.Range("J86").Value = 6571.8335
$ws.Range("K86").Value = 6217.3335
$ws.Range("L86").Value = 6571.8335
$ws.Range("M86").Value = -5094.3335
$ws.Range("N86").Value = -8817.833500000001
$ws.Range("H89").Value = 6394.5835
$ws.Range("I89").Value = 6217.3335
$ws.Range("J89").Value = 6571.8335
$ws.Range("K89").Value = 31086.6675
$ws.Range("L89").Value = 32859.1675
$ws.Range("M89").Value = -25470.6675
$ws.Range("N89").Value = -44091.1675

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 255.81818
$ws.Range("I7").Value = 220.5
$ws.Range("K7").Value = 220.5
$ws.Range("M7").Value = -107.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 122699.07
$ws.Range("I5").Value = 553.28125
$ws.Range("J5").Value = 556995.25
$ws.Range("K5").Value = 1659.84375
$ws.Range("L5").Value = 1670985.75
$ws.Range("M5").Value = -1547.84375
$ws.Range("N5").Value = -1671209.75
$ws.Range("H34").Value = 5049
$ws.Range("J34").Value = 7549.5
$ws.Range("L34").Value = 22648.5
$ws.Range("N34").Value = -22816.5
$ws.Range("H39").Value = 4000
$ws.Range("J39").Value = 4000
$ws.Range("L39").Value = 12000
$ws.Range("N39").Value = -12588
$ws.Range("H55").Value = 27375
$ws.Range("J55").Value = 34166.668
$ws.Range("L55").Value = 102500.004
$ws.Range("N55").Value = -102854.004
$ws.Range("H135").Value = 122699.07
$ws.Range("I135").Value = 553.28125
$ws.Range("J135").Value = 556995.25
$ws.Range("K135").Value = 4979.53125
$ws.Range("L135").Value = 5012957.25
$ws.Range("M135").Value = -2444.53125
$ws.Range("N135").Value = -5018027.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 17434.5
$ws.Range("I80").Value = 5000
$ws.Range("J80").Value = 21579.334
$ws.Range("K80").Value = 5000
$ws.Range("L80").Value = 21579.334
$ws.Range("M80").Value = -4002
$ws.Range("N80").Value = -23575.334
$ws.Range("H83").Value = 17434.5
$ws.Range("I83").Value = 5000
$ws.Range("J83").Value = 21579.334
$ws.Range("K83").Value = 25000
$ws.Range("L83").Value = 107896.67
$ws.Range("M83").Value = -20008
$ws.Range("N83").Value = -117880.67

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1203.84
$ws.Range("I93").Value = 1167.0769
$ws.Range("J93").Value = 1243.6666
$ws.Range("K93").Value = 1167.0769
$ws.Range("L93").Value = 1243.6666
$ws.Range("M93").Value = 80.92309999999998
$ws.Range("N93").Value = -3739.6666

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10816.728
$ws.Range("J62").Value = 10850.333
$ws.Range("L62").Value = 10850.333
$ws.Range("N62").Value = -12098.333
$ws.Range("H65").Value = 10816.728
$ws.Range("J65").Value = 10850.333
$ws.Range("L65").Value = 54251.665
$ws.Range("N65").Value = -60491.665
$ws.Range("H81").Value = 1500
$ws.Range("J81").Value = 1500
$ws.Range("L81").Value = 3000
$ws.Range("N81").Value = -5122
$ws.Range("H84").Value = 1500
$ws.Range("J84").Value = 1500
$ws.Range("L84").Value = 15000
$ws.Range("N84").Value = -25608
$ws.Range("H95").Value = 28822
$ws.Range("J95").Value = 28822
$ws.Range("L95").Value = 28822
$ws.Range("N95").Value = -34314
